$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 20 and 21: the Id (A), Lokalnamn (P), Ost (Q) and Nord (R) values were
# swapped between the two observation records.
foreach ($col in 1, 16, 17, 18) {
    $v1 = $ws.Cells.Item(20, $col).Value2
    $v2 = $ws.Cells.Item(21, $col).Value2
    $ws.Cells.Item(20, $col).Value2 = $v2
    $ws.Cells.Item(21, $col).Value2 = $v1
}

# Rows 33 and 34: the whole species records (Id, Taxonsorteringsordning,
# Rodlistade, TaxonId, Artnamn, Vetenskapligt namn, Auktor, Lokalnamn, Ost,
# Nord) were swapped between the two rows.
foreach ($col in 1, 2, 4, 5, 6, 7, 8, 16, 17, 18) {
    $v1 = $ws.Cells.Item(33, $col).Value2
    $v2 = $ws.Cells.Item(34, $col).Value2
    $ws.Cells.Item(33, $col).Value2 = $v2
    $ws.Cells.Item(34, $col).Value2 = $v1
}

# The Enhet (J), Kon (L), Aktivitet (M), Metod (N) and Bestamningsmetod (AF)
# fields only applied to the barkbock record, which is now on row 33.
$ws.Cells.Item(34, 10).ClearContents()
$ws.Cells.Item(34, 12).ClearContents()
$ws.Cells.Item(34, 13).ClearContents()
$ws.Cells.Item(34, 14).ClearContents()
$ws.Cells.Item(34, 32).ClearContents()

$ws.Cells.Item(33, 13).Value2 = "äldre gnagspår"

Write-Host "Row swap complete"
